$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the "799" block, pushing all existing
# data (rows 799-888) down to rows 802-891.
$ws.Range("A799:A801").EntireRow.Insert()

# --- Row 799: new record, 1a amarillo ---
$ws.Range("A799").Value = 2
$ws.Range("B799").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C799").Value = "Coquimbo"
$ws.Range("D799").Value = 45106
$ws.Range("E799").Value = 4
$ws.Range("F799").Value = "Fruta"
$ws.Range("G799").Value = 100102
$ws.Range("H799").Value = "Cítricos"
$ws.Range("I799").Value = 100102003
$ws.Range("J799").Value = "Limón"
$ws.Range("K799").Value = "Sin especificar"
$ws.Range("L799").Value = "1a amarillo"
$ws.Range("M799").Value = 560
$ws.Range("N799").Value = 5300
$ws.Range("O799").Value = 5500
$ws.Range("P799").Value = 5400
$ws.Range("Q799").Value = "$/malla 18 kilos"
$ws.Range("R799").Value = "Provincia de Limarí"
$ws.Range("S799").Value = 300
$ws.Range("T799").Value = 18

# --- Row 800: new record, 2a amarillo ---
$ws.Range("A800").Value = 2
$ws.Range("B800").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C800").Value = "Coquimbo"
$ws.Range("D800").Value = 45106
$ws.Range("E800").Value = 4
$ws.Range("F800").Value = "Fruta"
$ws.Range("G800").Value = 100102
$ws.Range("H800").Value = "Cítricos"
$ws.Range("I800").Value = 100102003
$ws.Range("J800").Value = "Limón"
$ws.Range("K800").Value = "Sin especificar"
$ws.Range("L800").Value = "2a amarillo"
$ws.Range("M800").Value = 420
$ws.Range("N800").Value = 3300
$ws.Range("O800").Value = 3500
$ws.Range("P800").Value = 3400
$ws.Range("Q800").Value = "$/malla 18 kilos"
$ws.Range("R800").Value = "Provincia de Limarí"
$ws.Range("S800").Value = 189
$ws.Range("T800").Value = 18

# --- Row 801: new record, 3a amarillo ---
$ws.Range("A801").Value = 2
$ws.Range("B801").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C801").Value = "Coquimbo"
$ws.Range("D801").Value = 45106
$ws.Range("E801").Value = 4
$ws.Range("F801").Value = "Fruta"
$ws.Range("G801").Value = 100102
$ws.Range("H801").Value = "Cítricos"
$ws.Range("I801").Value = 100102003
$ws.Range("J801").Value = "Limón"
$ws.Range("K801").Value = "Sin especificar"
$ws.Range("L801").Value = "3a amarillo"
$ws.Range("M801").Value = 320
$ws.Range("N801").Value = 1300
$ws.Range("O801").Value = 1500
$ws.Range("P801").Value = 1400
$ws.Range("Q801").Value = "$/malla 18 kilos"
$ws.Range("R801").Value = "Provincia de Limarí"
$ws.Range("S801").Value = 78
$ws.Range("T801").Value = 18

# Match the date-formatted style used by every other "Fecha" cell in
# column D (row 802 is the first untouched original row after the insert).
$ws.Range("D799:D801").NumberFormat = $ws.Cells.Item(802, 4).NumberFormat
